$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Formula = "=AVERAGE(B2:B5)"
$ws.Range("C6").Formula = "=B6/100"

$ws.Range("B6:C6").Select()
